# Fix "Recorded By" column G order: swap "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
